$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for Matteo's effort table (before the "Total effort" row)
# to record the "Use case diagrams" work item.
$ws.Rows("10:10").Insert()

# Copy formatting from the row above (same table) so the new row matches
# the existing striped/bordered style used by the other data rows.
$ws.Range("A9:C9").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122)

# Fill in the new row's data.
$ws.Range("A10").Value = 43770
$ws.Range("B10").Value = "Use case diagrams"
$ws.Range("C10").Value = 2

# Extend the "Total effort" SUM formula to include the newly inserted row.
$ws.Range("C11").Formula = "=SUM(C4:C10)"

# Restore the selection to match the latest save state.
$ws.Range("F9").Select() | Out-Null
